$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numRows = 24
$bd = New-Object "object[,]" $numRows,3
$fk = New-Object "object[,]" $numRows,6
$m  = New-Object "object[,]" $numRows,1

$bd[0,0] = 0.3128181661851386; $bd[0,1] = 0.06058390144613668; $bd[0,2] = 0.1854047359384765
$fk[0,0] = 2.017967820664794; $fk[0,1] = 1.280923242651482; $fk[0,2] = 1.218535219912155; $fk[0,3] = 1.028265946199681; $fk[0,4] = 0.3007937144052946; $fk[0,5] = 0.3584682315354826
$m[0,0] = 0.2880942091226117
$bd[1,0] = 0.2837422116946868; $bd[1,1] = 0.05482136769936119; $bd[1,2] = 0.1820666772606785
$fk[1,0] = 2.016264295554791; $fk[1,1] = 1.279475906418924; $fk[1,2] = 1.222637093007506; $fk[1,3] = 1.032030271367667; $fk[1,4] = 0.2988999196909603; $fk[1,5] = 0.3249661433025324
$m[1,0] = 0.2761836987766415
$bd[2,0] = 0.2659796415130984; $bd[2,1] = 0.05130270292833927; $bd[2,2] = 0.1800945990049598
$fk[2,0] = 2.016183255596644; $fk[2,1] = 1.279274012295289; $fk[2,2] = 1.225641824270426; $fk[2,3] = 1.034805479855429; $fk[2,4] = 0.2978817543344618; $fk[2,5] = 0.3045019095141868
$m[2,0] = 0.2690126815871352
$bd[3,0] = 0.2587642414986249; $bd[3,1] = 0.04987375105439185; $bd[3,2] = 0.1793105289186201
$fk[3,0] = 2.01639287905553; $fk[3,1] = 1.279364293606179; $fk[3,2] = 1.226988570100971; $fk[3,3] = 1.036053010757829; $fk[3,4] = 0.2975032243246432; $fk[3,5] = 0.2961895697887371
$m[3,0] = 0.2661262823297861
$bd[4,0] = 0.2575675277689982; $bd[4,1] = 0.04963677307306114; $bd[4,2] = 0.1791815188591386
$fk[4,0] = 2.01644234249008; $fk[4,1] = 1.279389702536164; $fk[4,2] = 1.227219583737138; $fk[4,3] = 1.03626720416845; $fk[4,4] = 0.2974425673611307; $fk[4,5] = 0.2948109518805353
$m[4,0] = 0.2656491666586049
$bd[5,0] = 0.2658822385944859; $bd[5,1] = 0.05128341159002048; $bd[5,2] = 0.180083945427171
$fk[5,0] = 2.01618510017353; $fk[5,1] = 1.279274531382853; $fk[5,2] = 1.225659491740089; $fk[5,3] = 1.034821832396318; $fk[5,4] = 0.2978765020266394; $fk[5,5] = 0.3043896966633355
$m[5,0] = 0.2689736092875918
$bd[6,0] = 0.3027742379424581; $bd[6,1] = 0.05859292865774535; $bd[6,2] = 0.1842377333313863
$fk[6,0] = 2.017180172852989; $fk[6,1] = 1.280281572367556; $fk[6,2] = 1.219848666833101; $fk[6,3] = 1.029467578288035; $fk[6,4] = 0.300110721193434; $fk[6,5] = 0.3468948420494371
$m[6,0] = 0.2839580430582913
$bd[7,0] = 0.3758249848788466; $bd[7,1] = 0.07308212949126869; $bd[7,2] = 0.1929953883095124
$fk[7,0] = 2.026790075390124; $fk[7,1] = 1.287713742595813; $fk[7,2] = 1.212309958029465; $fk[7,3] = 1.022651176341036; $fk[7,4] = 0.3056396971644375; $fk[7,5] = 0.4310810987830109
$m[7,0] = 0.3144664128324948
$bd[8,0] = 0.4299179615258311; $bd[8,1] = 0.08382330431616936; $bd[8,2] = 0.1997998284582252
$fk[8,0] = 2.038526293476565; $fk[8,1] = 1.296514878304492; $fk[8,2] = 1.209121621695772; $fk[8,3] = 1.019892934790043; $fk[8,4] = 0.3104024993508006; $fk[8,5] = 0.493435959060804
$m[8,0] = 0.3375642107448442
$bd[9,0] = 0.4546167971069792; $bd[9,1] = 0.08873093677581778; $bd[9,2] = 0.2029751637502812
$fk[9,0] = 2.044882475226544; $fk[9,1] = 1.301247276331011; $fk[9,2] = 1.208181498654284; $fk[9,3] = 1.019127664649602; $fk[9,4] = 0.312721640864865; $fk[9,5] = 0.5219114371939781
$m[9,0] = 0.3482201529563582
$bd[10,0] = 0.4639825570927201; $bd[10,1] = 0.09059241339443247; $bd[10,2] = 0.2041890135200788
$fk[10,0] = 2.04743575857826; $fk[10,1] = 1.303144306255817; $fk[10,2] = 1.207898858329457; $fk[10,3] = 1.018908330082297; $fk[10,4] = 0.3136217766800939; $fk[10,5] = 0.5327099695520303
$m[10,0] = 0.3522765776480483
$bd[11,0] = 0.461964907041164; $bd[11,1] = 0.0901913751523864; $bd[11,2] = 0.2039270826372643
$fk[11,0] = 2.046879354282737; $fk[11,1] = 1.30273107525359; $fk[11,2] = 1.207956467285015; $fk[11,3] = 1.018952432983845; $fk[11,4] = 0.3134269411808219; $fk[11,5] = 0.5303836287269519
$m[11,0] = 0.3514020108401041
$bd[12,0] = 0.4553870682557033; $bd[12,1] = 0.08888402029367626; $bd[12,2] = 0.2030747994955249
$fk[12,0] = 2.045089602341946; $fk[12,1] = 1.301401241241848; $fk[12,2] = 1.208156775409378; $fk[12,3] = 1.01910820739576; $fk[12,4] = 0.3127952561522278; $fk[12,5] = 0.5227995304101682
$m[12,0] = 0.3485534522663798
$bd[13,0] = 0.4513596152226285; $bd[13,1] = 0.08808362604537479; $bd[13,2] = 0.202554236238754
$fk[13,0] = 2.044012386149745; $fk[13,1] = 1.300600355640555; $fk[13,2] = 1.208289023717413; $fk[13,3] = 1.019212801384683; $fk[13,4] = 0.3124111862693297; $fk[13,5] = 0.5181560585142222
$m[13,0] = 0.3468113917024525
$bd[14,0] = 0.4283056504710885; $bd[14,1] = 0.08350300850220549; $bd[14,2] = 0.1995939150814934
$fk[14,0] = 2.038131377025522; $fk[14,1] = 1.296220284425701; $fk[14,2] = 1.209193327791979; $fk[14,3] = 1.019952802467692; $fk[14,4] = 0.3102540066575727; $fk[14,5] = 0.4915772042483013
$m[14,0] = 0.3368708002793355
$bd[15,0] = 0.4141860294763262; $bd[15,1] = 0.08069842135303418; $bd[15,2] = 0.1977982775140816
$fk[15,0] = 2.034784160630522; $fk[15,1] = 1.293720026200177; $fk[15,2] = 1.209878772026627; $fk[15,3] = 1.020532184629708; $fk[15,4] = 0.308969707143433; $fk[15,5] = 0.4752998660499088
$m[15,0] = 0.3308105512066959
$bd[16,0] = 0.4060734451499286; $bd[16,1] = 0.07908731398465818; $bd[16,2] = 0.1967730053243741
$fk[16,0] = 2.032954682208668; $fk[16,1] = 1.29235052307277; $fk[16,2] = 1.210321050838772; $fk[16,3] = 1.020911496641794; $fk[16,4] = 0.3082453681115993; $fk[16,5] = 0.4659479387811132
$m[16,0] = 0.3273388540767854
$bd[17,0] = 0.4033281590718047; $bd[17,1] = 0.07854216766048694; $bd[17,2] = 0.1964271616585336
$fk[17,0] = 2.032351696334914; $fk[17,1] = 1.29189860576929; $fk[17,2] = 1.210479048323478; $fk[17,3] = 1.021047834994235; $fk[17,4] = 0.3080025851793664; $fk[17,5] = 0.4627833273233364
$m[17,0] = 0.3261658049942966
$bd[18,0] = 0.4156881946108228; $bd[18,1] = 0.08099676581828419; $bd[18,2] = 0.1979886473383203
$fk[18,0] = 2.035130567226616; $fk[18,1] = 1.293979083915488; $fk[18,2] = 1.209800834537674; $fk[18,3] = 1.020465740154521; $fk[18,4] = 0.3091049372542329; $fk[18,5] = 0.4770315458719381
$m[18,0] = 0.3314542272378702
$bd[19,0] = 0.4573187929510425; $bd[19,1] = 0.08926793912095832; $bd[19,2] = 0.203324826428485
$fk[19,0] = 2.045611324003502; $fk[19,1] = 1.301788995304477; $fk[19,2] = 1.208095949075158; $fk[19,3] = 1.019060539929747; $fk[19,4] = 0.3129802022955772; $fk[19,5] = 0.5250267461775593
$m[19,0] = 0.3493895672126754
$bd[20,0] = 0.4846015237491201; $bd[20,1] = 0.09469147315991222; $bd[20,2] = 0.2068788513962119
$fk[20,0] = 2.053314080068077; $fk[20,1] = 1.307505175344076; $fk[20,2] = 1.207409323566282; $fk[20,3] = 1.018552868878707; $fk[20,4] = 0.3156407187375834; $fk[20,5] = 0.5564845050051588
$m[20,0] = 0.3612351569769032
$bd[21,0] = 0.4700334893163927; $bd[21,1] = 0.09179520380857298; $bd[21,2] = 0.2049759405205123
$fk[21,0] = 2.049124911501053; $fk[21,1] = 1.304398287881526; $fk[21,2] = 1.207736665435831; $fk[21,3] = 1.018786218879789; $fk[21,4] = 0.3142090582574184; $fk[21,5] = 0.5396867669148548
$m[21,0] = 0.3549016548776365
$bd[22,0] = 0.415009050236705; $bd[22,1] = 0.08086188027567687; $bd[22,2] = 0.1979025591328991
$fk[22,0] = 2.034973661253176; $fk[22,1] = 1.293861752338088; $fk[22,2] = 1.209835919904805; $fk[22,3] = 1.020495635737582; $fk[22,4] = 0.3090437560449999; $fk[22,5] = 0.476248634304369
$m[22,0] = 0.3311631826068009
$bd[23,0] = 0.3559880713555685; $bd[23,1] = 0.06914566144595824; $bd[23,2] = 0.1905609893543811
$fk[23,0] = 2.023369722581435; $fk[23,1] = 1.285117471760842; $fk[23,2] = 1.213936594616726; $fk[23,3] = 1.024100425292517; $fk[23,4] = 0.3040209922058636; $fk[23,5] = 0.4310810987830109
$m[23,0] = 0.3060929628873623

$ws.Range("B2:D25").Value2 = $bd
$ws.Range("F2:K25").Value2 = $fk
$ws.Range("M2:M25").Value2 = $m
